$d = $word.ActiveDocument

# 1. Merge the three IMAGE runs into a single run with the same text.
$d.Content.Find.Execute(
    "***IMAGE imageGenerator(ara.logoUrl, 6, 2)***",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "***IMAGE imageGenerator(ara.logoUrl, 6, 2)***", 2
) | Out-Null

# 2. Fix typo "Jura" -> "Juros"
$d.Content.Find.Execute(
    "Jura:",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Juros:", 2
) | Out-Null
